$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress value for row 8 (SVM -> Unet -> Mask_RCNN task) and add comment
$ws.Range("C8").Value = 0.6
$ws.Range("D8").Value = "SVM -> Unet -> Mask_RCNN"

# Update progress value for row 13
$ws.Range("C13").Value = 0.25

# Update the selected / visible cell in the sheet view
$ws.Range("C14").Select()
$excel.ActiveWindow.ScrollRow = 7
